# UC3.1_TC1.xlsx - "test code generation module - update evaluations"
#
# Re-runs the CodeBLEU evaluation for this test case:
#  - the "Assertion validity" note (C7) is cleared (no longer flagged as
#    "not passing correctly as the baseline"), and the cell is re-styled
#    with an underlined font,
#  - the Code BLEU score (B12) and its breakdown dict (C12) are refreshed
#    with the new dataflow_match_score,
#  - the sheet is left with C7 selected and the print setup normalized to
#    portrait / A4.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("QuantitativeMetrics")

# --- Assertion validity note (C7): drop the stale note text ---------------
# Clearing keeps the cell numeric/empty; re-applying the font (underline)
# mints the new "note" style (border + wrap text + underlined font) that
# replaces the old plain note style.
$ws.Range("C7").ClearContents()
$ws.Range("C7").Font.Name = "Aptos Narrow"
$ws.Range("C7").Font.Size = 11
$ws.Range("C7").Font.Underline = $true

# --- Code BLEU metrics (B12/C12): refreshed evaluation numbers ------------
$ws.Range("B12").Value = 0.3040507380433932
$ws.Range("C12").Value = "{'codebleu': 0.3040507380433932, 'ngram_match_score': 0.10122066127141109, 'weighted_ngram_match_score': 0.1314325073523782, 'syntax_match_score': 0.7359307359307359, 'dataflow_match_score': 0.24761904761904763}"

# --- Page setup: portrait / A4 --------------------------------------------
$ws.PageSetup.Orientation = 1
$ws.PageSetup.PaperSize = 9

# --- Selection left on C7 ---------------------------------------------------
$ws.Activate()
$ws.Range("C7").Select()
